{"js": "// Word JS API (Office.js) script.\n// Body of: async (context) => { ... }\n//\n// Target edit: the trailing \"\u2026\" paragraph (which carries the\n// `_GoBack` bookmark) is split into three paragraphs:\n//   1) a duplicate \"\u2026\" paragraph (no bookmark), now tagged en-US\n//   2) a new Heading-2 paragraph \"\u0410\u043b\u0433\u043e\u0440\u0438\u0442\u043c \u0425\u0430\u0431\u0435\u0440\u043c\u0430\u043d\u0430\" / English subtitle\n//   3) the original \"\u2026\" + _GoBack bookmark paragraph, now tagged en-US\n//\n// We locate the last paragraph of the body (the \"\u2026\" one), insert an\n// empty placeholder paragraph right before it, and replace that\n// placeholder's OOXML with the two new paragraphs built above. Then we\n// replace the final paragraph's OOXML in place (re-adding its\n// bookmark ourselves) so it picks up the `en-US` language tagging on\n// both the paragraph mark and the run.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// OOXML for the two brand-new paragraphs inserted before the\n// existing \"\u2026\" / _GoBack paragraph.\nconst newParagraphsOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>\\u2026</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p>' +\n  '<w:pPr><w:pStyle w:val=\"2\"/><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n  '<w:r><w:t>\\u0410</w:t></w:r>' +\n  '<w:r><w:t>\\u043B\\u0433\\u043E\\u0440\\u0438\\u0442\\u043C</w:t></w:r>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>\\u0425\\u0430\\u0431\\u0435\\u0440\\u043C\\u0430\\u043D\\u0430</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:br/></w:r>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">deadlocks </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>avoidence</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>haberman\\'s</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> algorithm</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\n// OOXML for the final paragraph: same \"\u2026\" text as before plus the\n// original _GoBack bookmark, now with en-US language tagging on both\n// the paragraph mark and the run.\nconst lastParagraphOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>\\u2026</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"1\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"1\"/>' +\n  '</w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\n// 1) Insert an empty placeholder paragraph before the last (\"\u2026\")\n//    paragraph, then replace the placeholder's content with the two\n//    new paragraphs (duplicate \"\u2026\" + the new heading).\nconst placeholder = lastParagraph.insertParagraph(\"\", Word.InsertLocation.before);\nawait context.sync();\n\nplaceholder.insertOoxml(newParagraphsOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Replace the original trailing paragraph in place so it keeps its\n//    bookmark but gains the en-US language tagging. Re-query the\n//    paragraphs collection first: the earlier sibling insertions above\n//    mean the original `lastParagraph` reference can no longer be\n//    trusted to resolve to the true last paragraph.\nparagraphs.load(\"items\");\nawait context.sync();\nconst trailingParagraph = paragraphs.items[paragraphs.items.length - 1];\n\ntrailingParagraph.insertOoxml(lastParagraphOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n#\n# Target edit: the trailing \"\u2026\" paragraph (which carries the\n# `_GoBack` bookmark) is split into three paragraphs:\n#   1) a duplicate \"\u2026\" paragraph (no bookmark), now tagged en-US\n#   2) a new Heading-2 paragraph \"\u0410\u043b\u0433\u043e\u0440\u0438\u0442\u043c \u0425\u0430\u0431\u0435\u0440\u043c\u0430\u043d\u0430\" / English subtitle\n#   3) the original \"\u2026\" + _GoBack bookmark paragraph, now tagged en-US\n#\n# We replace the Range of the last (\"\u2026\") paragraph with WordML for all\n# three paragraphs in a single InsertXML call (re-adding the _GoBack\n# bookmark ourselves, since it lived inside the paragraph being\n# replaced).\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$lastParagraph = $d.Paragraphs.Item($count)\n$targetRange = $lastParagraph.Range\n\n$wNs = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n\n# Paragraph 1: duplicate \"\u2026\" (en-US on paragraph mark + run).\n$para1 = '<w:p ' + $wNs + '>' + `\n    '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' + `\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>&#8230;</w:t></w:r>' + `\n  '</w:p>'\n\n# Paragraph 2: new Heading 2 \"\u0410\u043b\u0433\u043e\u0440\u0438\u0442\u043c \u0425\u0430\u0431\u0435\u0440\u043c\u0430\u043d\u0430\" + English subtitle.\n$para2 = '<w:p ' + $wNs + '>' + `\n    '<w:pPr><w:pStyle w:val=\"2\"/><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' + `\n    '<w:r><w:t>&#1040;</w:t></w:r>' + `\n    '<w:r><w:t>&#1083;&#1075;&#1086;&#1088;&#1080;&#1090;&#1084;</w:t></w:r>' + `\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>' + `\n    '<w:proofErr w:type=\"spellStart\"/>' + `\n    '<w:r><w:t>&#1061;&#1072;&#1073;&#1077;&#1088;&#1084;&#1072;&#1085;&#1072;</w:t></w:r>' + `\n    '<w:proofErr w:type=\"spellEnd\"/>' + `\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:br/></w:r>' + `\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">deadlocks </w:t></w:r>' + `\n    '<w:proofErr w:type=\"spellStart\"/>' + `\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>avoidence</w:t></w:r>' + `\n    '<w:proofErr w:type=\"spellEnd\"/>' + `\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>' + `\n    '<w:proofErr w:type=\"spellStart\"/>' + `\n    \"<w:r><w:rPr><w:lang w:val=`\"en-US`\"/></w:rPr><w:t>haberman's</w:t></w:r>\" + `\n    '<w:proofErr w:type=\"spellEnd\"/>' + `\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> algorithm</w:t></w:r>' + `\n  '</w:p>'\n\n# Paragraph 3: original \"\u2026\" text + the original _GoBack bookmark, now\n# tagged en-US on both the paragraph mark and the run.\n$para3 = '<w:p ' + $wNs + '>' + `\n    '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' + `\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>&#8230;</w:t></w:r>' + `\n    '<w:bookmarkStart w:id=\"1\" w:name=\"_GoBack\"/>' + `\n    '<w:bookmarkEnd w:id=\"1\"/>' + `\n  '</w:p>'\n\n$targetRange.InsertXML($para1 + $para2 + $para3)\n"}
